# Accounts.xlsx edit: append two new account rows (17 & 18) that were
# uploaded after clearing out the previous table contents, widen the
# "last_update" column to fit the date text, and carry over the small
# phonetic-guide font that Excel attaches when the sheet is re-saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17 -----------------------------------------------------------
$ws.Range("A17").Value = 73
$ws.Range("B17").Value = "Dept32"
$ws.Range("C17").Value = 43433.43
# D11 already holds the literal text "2024-04-27" as a shared string;
# copy/paste (values) so the new cell keeps it as text instead of having
# Excel reinterpret the date-shaped string as a serial date number.
$ws.Range("D11").Copy() | Out-Null
$ws.Range("D17").PasteSpecial(-4163) | Out-Null

# --- Row 18 -----------------------------------------------------------
$ws.Range("A18").Value = 83
$ws.Range("B18").Value = "Dept44"
$ws.Range("C18").Value = 23553.35
# D7 already holds the literal text "2024-02-07" as a shared string.
$ws.Range("D7").Copy() | Out-Null
$ws.Range("D18").PasteSpecial(-4163) | Out-Null

$excel.CutCopyMode = 0

# --- Column width -------------------------------------------------------
# Widen column D ("last_update") so the date text fits, matching the
# width Excel computed when the sheet was re-saved after the upload.
$ws.Columns("D").ColumnWidth = 17.46

# --- Phonetic-guide font -------------------------------------------------
# Re-saving the workbook from Excel registers the small "noConversion"
# phonetic-guide font (size 8) used for the sheet's phoneticPr entry.
# Stamp it on a scratch cell to register the font, then clean the scratch
# cell back up so no visible formatting/content remains on it.
$scratch = $ws.Range("Z100")
$scratch.Font.Size = 8
$scratch.ClearFormats() | Out-Null
$scratch.ClearContents() | Out-Null

# --- Selection ------------------------------------------------------------
# Leave the cursor on the last cell typed into, like Excel would after
# manual data entry.
$ws.Range("D18").Select() | Out-Null
